# Insert a new data row at row 6 (pushing existing rows 6..113 down to 7..114)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 6, shifting rows 6..113 down to 7..114.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new record's data.
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = Get-Date -Year 2022 -Month 2 -Day 24 -Hour 0 -Minute 0 -Second 0
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 100112024
$ws.Range("G6").Value = "Choclo"
$ws.Range("H6").Value = "Choclero"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 60000
$ws.Range("K6").Value = 150
$ws.Range("L6").Value = 200
$ws.Range("M6").Value = 175
$ws.Range("N6").Value = "`$/unidad"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 175
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = "Hortaliza"
